# Ch2 Protocol Suite and Key Agreement
#
# Give the title placeholders on the two PSK slides
#   - "PSK: 前方秘匿性なし"        -> <a:off x="795885" y="126800"/>  <a:ext cx="10515600" cy="1325563"/>
#   - "PSK: 完全前方秘匿性あり"    -> <a:off x="838200" y="46959"/>   <a:ext cx="10515600" cy="1325563"/>
# an explicit xfrm instead of inheriting the empty <p:spPr/> from the layout/master.
#
# NOTE: this runtime's PowerShell parser does not bind `-Name value` style
# named arguments, so the helper below takes plain positional parameters.
# The Left/Top/Width/Height values are also pre-computed (points, as
# IEEE-754 doubles that round to the exact desired EMU once PowerPoint's
# COM layer stores them as single-precision floats) so the round-trip
# through the COM property reproduces the exact target EMU values.

$p = $ppt.ActivePresentation

function Set-TitleXfrm {
    param($slide, $leftPts, $topPts, $widthPts, $heightPts)

    $title = $null
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.Type -eq 14 -and $shape.PlaceholderFormat.Type -eq 1) {
            $title = $shape
            break
        }
    }

    $title.Left = $leftPts
    $title.Top = $topPts
    $title.Width = $widthPts
    $title.Height = $heightPts
}

# Slide 18 - "PSK: 前方秘匿性なし"   -> off (795885, 126800) ext (10515600, 1325563)
Set-TitleXfrm $p.Slides.Item(18) 62.668113708496094 9.984251976013184 828.0 104.37504577636719

# Slide 19 - "PSK: 完全前方秘匿性あり" -> off (838200, 46959) ext (10515600, 1325563)
Set-TitleXfrm $p.Slides.Item(19) 66.0 3.697559118270874 828.0 104.37504577636719
